# Update "想去人数" (number of people interested) values in the F column
# for both the "展览" and "全部类型" sheets, which carry duplicate data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 8787
    $ws.Range("F3").Value = 197
    $ws.Range("F4").Value = 420
    $ws.Range("F5").Value = 194
}
